# Data driven Testing for LoginTest
#
# The Login sheet previously held two distinct shared strings ("admin" in
# A1 and "admin123" in B1). For data-driven testing both cells are now
# populated with the same test value ("tutorial"), which also collapses
# the shared-strings table down to a single unique entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "tutorial"
$ws.Range("B1").Value = "tutorial"

# Move the sheet's active selection from D5 to D6.
$ws.Range("D6").Select()
